# Refresh cryptos list: updates Price (D) / Volume(1h) (E) figures for
# every coin row, plus a 3-way reshuffle of rows 45-47 (new entrant
# EnergySwap pushes VeChain and FraxShare down one spot).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of Price cells are plain decimals (e.g. "1.00", "18.23").
# Mark them as Text first so Excel keeps the exact string instead of
# silently reinterpreting them as numbers and dropping trailing zeros.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D12", "D16", "D18", "D20", "D23", "D24", "D25", "D26", "D28", "D30", "D31", "D32", "D34", "D35", "D37", "D38", "D39", "D44", "D45", "D46", "D47", "D48", "D50")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "40.170.10"
$ws.Range("E2").Value = "  +1.04%  "

# Row 3
$ws.Range("D3").Value = "2.236.60"
$ws.Range("E3").Value = "  -0.09%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "293.71"
$ws.Range("E5").Value = "  -1.51%  "

# Row 6
$ws.Range("D6").Value = "88.67"
$ws.Range("E6").Value = "  +5.92%  "

# Row 7
$ws.Range("D7").Value = "0.519"
$ws.Range("E7").Value = "  +0.26%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").Value = "0.477"
$ws.Range("E9").Value = "  +0.22%  "

# Row 10
$ws.Range("D10").Value = "31.19"
$ws.Range("E10").Value = "  +2.43%  "

# Row 11
$ws.Range("E11").Value = "  +1.26%  "

# Row 12
$ws.Range("D12").Value = "47.72"
$ws.Range("E12").Value = "  +1.19%  "

# Row 13
$ws.Range("E13").Value = "  +1.56%  "

# Row 14
$ws.Range("E14").Value = "  +2.00%  "

# Row 15
$ws.Range("D15").Value = "2.579.91"
$ws.Range("E15").Value = "  -0.35%  "

# Row 16
$ws.Range("D16").Value = "14.22"
$ws.Range("E16").Value = "  -0.21%  "

# Row 17
$ws.Range("D17").Value = "2.259.84"
$ws.Range("E17").Value = "  +0.80%  "

# Row 18
$ws.Range("D18").Value = "0.738"
$ws.Range("E18").Value = "  +2.16%  "

# Row 19
$ws.Range("D19").Value = "40.122.21"
$ws.Range("E19").Value = "  +1.05%  "

# Row 20
$ws.Range("D20").Value = "11.66"
$ws.Range("E20").Value = "  +10.92%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0891"
$ws.Range("E21").Value = "  +0.94%  "

# Row 22
$ws.Range("E22").Value = "  +0.84%  "

# Row 23
$ws.Range("D23").Value = "66.23"
$ws.Range("E23").Value = "  +1.23%  "

# Row 24
$ws.Range("D24").Value = "236.67"
$ws.Range("E24").Value = "  +3.28%  "

# Row 25
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.12%  "

# Row 26
$ws.Range("D26").Value = "2.48"
$ws.Range("E26").Value = "  +2.17%  "

# Row 28
$ws.Range("D28").Value = "23.03"
$ws.Range("E28").Value = "  +0.24%  "

# Row 29
$ws.Range("E29").Value = "  +3.80%  "

# Row 30
$ws.Range("D30").Value = "9.35"
$ws.Range("E30").Value = "  +1.63%  "

# Row 31
$ws.Range("D31").Value = "33.15"
$ws.Range("E31").Value = "  +1.15%  "

# Row 32
$ws.Range("D32").Value = "152.64"
$ws.Range("E32").Value = "  +1.76%  "

# Row 33
$ws.Range("E33").Value = "  -0.21%  "

# Row 34
$ws.Range("D34").Value = "5.00"
$ws.Range("E34").Value = "  +2.36%  "

# Row 35
$ws.Range("D35").Value = "0.0724"
$ws.Range("E35").Value = "  +2.49%  "

# Row 36
$ws.Range("E36").Value = "  -2.24%  "

# Row 37
$ws.Range("D37").Value = "2.88"
$ws.Range("E37").Value = "  +7.37%  "

# Row 38
$ws.Range("D38").Value = "16.27"
$ws.Range("E38").Value = "  +0.06%  "

# Row 39
$ws.Range("D39").Value = "0.112"
$ws.Range("E39").Value = "  +0.95%  "

# Row 40
$ws.Range("E40").Value = "  +2.98%  "

# Row 41
$ws.Range("E41").Value = "  +3.01%  "

# Row 42
$ws.Range("D42").Value = "2.112.55"
$ws.Range("E42").Value = "  +9.60%  "

# Row 43
$ws.Range("E43").Value = "  +2.84%  "

# Row 44
$ws.Range("D44").Value = "2.18"
$ws.Range("E44").Value = "  +6.57%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "18.23"
$ws.Range("E45").Value = "  +9.89%  "

# Row 46
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0272"
$ws.Range("E46").Value = "  +2.68%  "

# Row 47
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "10.11"
$ws.Range("E47").Value = "  +10.41%  "

# Row 48
$ws.Range("D48").Value = "2.69"
$ws.Range("E48").Value = "  +2.08%  "

# Row 49
$ws.Range("D49").Value = "2.448.67"
$ws.Range("E49").Value = "  -0.49%  "

# Row 50
$ws.Range("D50").Value = "71.49"
$ws.Range("E50").Value = "  -0.26%  "

# Row 51
$ws.Range("E51").Value = "  +6.03%  "
